$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(92, 8).Value2 = 1770
$ws.Cells.Item(92, 10).Value2 = 0
$ws.Cells.Item(92, 12).Value2 = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(103, 8).Value2 = 633.375
$ws.Cells.Item(103, 10).Value2 = 695.1667
$ws.Cells.Item(103, 12).Value2 = 2085.5001
$ws.Cells.Item(103, 14).Value2 = -3257.5001
$ws.Cells.Item(107, 8).Value2 = 3636.25
$ws.Cells.Item(107, 10).Value2 = 5141.2
$ws.Cells.Item(107, 12).Value2 = 5141.2
$ws.Cells.Item(107, 14).Value2 = -8981.200000000001
$ws.Cells.Item(135, 8).Value2 = 791.1429000000001
$ws.Cells.Item(135, 9).Value2 = 750.3333
$ws.Cells.Item(135, 11).Value2 = 6752.9997
$ws.Cells.Item(135, 13).Value2 = -4217.9997
$ws.Cells.Item(137, 8).Value2 = 3148.375
$ws.Cells.Item(137, 9).Value2 = 3299
$ws.Cells.Item(137, 10).Value2 = 3058
$ws.Cells.Item(137, 11).Value2 = 9897
$ws.Cells.Item(137, 12).Value2 = 9174
$ws.Cells.Item(137, 13).Value2 = -7347
$ws.Cells.Item(137, 14).Value2 = -14274

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value2 = 2740
$ws.Cells.Item(2, 9).Value2 = 2038
$ws.Cells.Item(2, 10).Value2 = 4495
$ws.Cells.Item(2, 11).Value2 = 2038
$ws.Cells.Item(2, 12).Value2 = 4495
$ws.Cells.Item(2, 13).Value2 = -1925
$ws.Cells.Item(2, 14).Value2 = -4721
$ws.Cells.Item(37, 8).Value2 = 21250
$ws.Cells.Item(37, 9).Value2 = 10000
$ws.Cells.Item(37, 10).Value2 = 25000
$ws.Cells.Item(37, 11).Value2 = 10000
$ws.Cells.Item(37, 12).Value2 = 25000
$ws.Cells.Item(37, 13).Value2 = -9727
$ws.Cells.Item(37, 14).Value2 = -25546
$ws.Cells.Item(61, 8).Value2 = 2488.125
$ws.Cells.Item(61, 9).Value2 = 1772.1428
$ws.Cells.Item(61, 11).Value2 = 1772.1428
$ws.Cells.Item(61, 13).Value2 = -1560.1428
$ws.Cells.Item(74, 8).Value2 = 2904
$ws.Cells.Item(74, 9).Value2 = 6712
$ws.Cells.Item(74, 10).Value2 = 1000
$ws.Cells.Item(74, 11).Value2 = 6712
$ws.Cells.Item(74, 12).Value2 = 1000
$ws.Cells.Item(74, 13).Value2 = -5838
$ws.Cells.Item(74, 14).Value2 = -2748
$ws.Cells.Item(77, 8).Value2 = 2904
$ws.Cells.Item(77, 9).Value2 = 6712
$ws.Cells.Item(77, 10).Value2 = 1000
$ws.Cells.Item(77, 11).Value2 = 33560
$ws.Cells.Item(77, 12).Value2 = 5000
$ws.Cells.Item(77, 13).Value2 = -29192
$ws.Cells.Item(77, 14).Value2 = -13736
$ws.Cells.Item(110, 8).Value2 = 2787.8
$ws.Cells.Item(110, 9).Value2 = 516.25
$ws.Cells.Item(110, 11).Value2 = 516.25
$ws.Cells.Item(110, 13).Value2 = 1528.75
$ws.Cells.Item(116, 8).Value2 = 2740
$ws.Cells.Item(116, 9).Value2 = 2038
$ws.Cells.Item(116, 10).Value2 = 4495
$ws.Cells.Item(116, 11).Value2 = 2038
$ws.Cells.Item(116, 12).Value2 = 4495
$ws.Cells.Item(116, 13).Value2 = 256
$ws.Cells.Item(116, 14).Value2 = -9083
$ws.Cells.Item(122, 8).Value2 = 1107
$ws.Cells.Item(122, 10).Value2 = 1107
$ws.Cells.Item(122, 12).Value2 = 3321
$ws.Cells.Item(122, 14).Value2 = -8221
$ws.Cells.Item(132, 8).Value2 = 894.125
$ws.Cells.Item(132, 9).Value2 = 894.125
$ws.Cells.Item(132, 11).Value2 = 2682.375
$ws.Cells.Item(132, 13).Value2 = -152.375
$ws.Cells.Item(136, 8).Value2 = 2488.125
$ws.Cells.Item(136, 9).Value2 = 1772.1428
$ws.Cells.Item(136, 11).Value2 = 5316.428400000001
$ws.Cells.Item(136, 13).Value2 = -2766.428400000001

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value2 = 2740
$ws.Cells.Item(3, 9).Value2 = 2038
$ws.Cells.Item(3, 10).Value2 = 4495
$ws.Cells.Item(3, 11).Value2 = 2038
$ws.Cells.Item(3, 12).Value2 = 4495
$ws.Cells.Item(3, 13).Value2 = -1924
$ws.Cells.Item(3, 14).Value2 = -4723
$ws.Cells.Item(16, 8).Value2 = 204.28572
$ws.Cells.Item(16, 9).Value2 = 0
$ws.Cells.Item(16, 10).Value2 = 204.28572
$ws.Cells.Item(16, 11).Value2 = 0
$ws.Cells.Item(16, 12).ClearContents()
$ws.Cells.Item(16, 13).Value2 = 204.28572
$ws.Cells.Item(16, 14).Value2 = -544.28572
$ws.Cells.Item(134, 8).Value2 = 6476.231
$ws.Cells.Item(134, 9).Value2 = 6476.231
$ws.Cells.Item(134, 11).Value2 = 19428.693
$ws.Cells.Item(134, 13).Value2 = -16893.693

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value2 = 2476.3333
$ws.Cells.Item(31, 9).Value2 = 2476.3333
$ws.Cells.Item(31, 11).Value2 = 2476.3333
$ws.Cells.Item(31, 13).Value2 = -2181.3333
$ws.Cells.Item(34, 8).Value2 = 2476.3333
$ws.Cells.Item(34, 9).Value2 = 2476.3333
$ws.Cells.Item(34, 11).Value2 = 2476.3333
$ws.Cells.Item(34, 13).Value2 = -2274.3333
$ws.Cells.Item(58, 8).Value2 = 5816
$ws.Cells.Item(58, 9).Value2 = 6285.3335
$ws.Cells.Item(58, 10).Value2 = 3000
$ws.Cells.Item(58, 11).Value2 = 6285.3335
$ws.Cells.Item(58, 12).Value2 = 3000
$ws.Cells.Item(58, 13).Value2 = -6082.3335
$ws.Cells.Item(58, 14).Value2 = -3406
$ws.Cells.Item(68, 8).Value2 = 40000
$ws.Cells.Item(68, 9).Value2 = 0
$ws.Cells.Item(68, 10).Value2 = 40000
$ws.Cells.Item(68, 11).Value2 = 0
$ws.Cells.Item(68, 12).ClearContents()
$ws.Cells.Item(68, 13).Value2 = 40000
$ws.Cells.Item(68, 14).Value2 = -41498
$ws.Cells.Item(71, 8).Value2 = 40000
$ws.Cells.Item(71, 9).Value2 = 0
$ws.Cells.Item(71, 10).Value2 = 40000
$ws.Cells.Item(71, 11).Value2 = 0
$ws.Cells.Item(71, 12).ClearContents()
$ws.Cells.Item(71, 13).Value2 = 120000
$ws.Cells.Item(71, 14).Value2 = -127488
$ws.Cells.Item(93, 8).Value2 = 34601.75
$ws.Cells.Item(93, 9).Value2 = 34601.75
$ws.Cells.Item(93, 11).Value2 = 34601.75
$ws.Cells.Item(93, 13).Value2 = -32729.75
$ws.Cells.Item(132, 8).Value2 = 2049
$ws.Cells.Item(132, 9).Value2 = 498
$ws.Cells.Item(132, 11).Value2 = 1494
$ws.Cells.Item(132, 13).Value2 = 1036
$ws.Cells.Item(134, 8).Value2 = 2927.6428
$ws.Cells.Item(134, 9).Value2 = 2927.6428
$ws.Cells.Item(134, 10).Value2 = 0
$ws.Cells.Item(134, 11).Value2 = 8782.928400000001
$ws.Cells.Item(134, 12).Value2 = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).Value2 = -6247.928400000001
$ws.Cells.Item(136, 8).Value2 = 5816
$ws.Cells.Item(136, 9).Value2 = 6285.3335
$ws.Cells.Item(136, 10).Value2 = 3000
$ws.Cells.Item(136, 11).Value2 = 18856.0005
$ws.Cells.Item(136, 12).Value2 = 9000
$ws.Cells.Item(136, 13).Value2 = -16306.0005
$ws.Cells.Item(136, 14).Value2 = -14100

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(39, 8).Value2 = 6666.6665
$ws.Cells.Item(39, 10).Value2 = 5000
$ws.Cells.Item(39, 12).Value2 = 15000
$ws.Cells.Item(39, 14).Value2 = -15588
$ws.Cells.Item(92, 8).Value2 = 840
$ws.Cells.Item(92, 10).Value2 = 750
$ws.Cells.Item(92, 12).Value2 = 2250
$ws.Cells.Item(92, 14).Value2 = -4746

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value2 = 4027
$ws.Cells.Item(132, 9).Value2 = 2081
$ws.Cells.Item(132, 11).Value2 = 6243
$ws.Cells.Item(132, 13).Value2 = -3713

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value2 = 1299.75
$ws.Cells.Item(16, 9).Value2 = 1342.5714
$ws.Cells.Item(16, 10).Value2 = 1000
$ws.Cells.Item(16, 11).Value2 = 1342.5714
$ws.Cells.Item(16, 12).Value2 = 1000
$ws.Cells.Item(16, 13).Value2 = -1172.5714
$ws.Cells.Item(16, 14).Value2 = -1340
$ws.Cells.Item(22, 8).Value2 = 2943.25
$ws.Cells.Item(22, 9).Value2 = 2918.75
$ws.Cells.Item(22, 10).Value2 = 2992.25
$ws.Cells.Item(22, 11).Value2 = 2918.75
$ws.Cells.Item(22, 12).Value2 = 2992.25
$ws.Cells.Item(22, 13).Value2 = -2623.75
$ws.Cells.Item(22, 14).Value2 = -3582.25
$ws.Cells.Item(24, 8).Value2 = 5650
$ws.Cells.Item(24, 9).Value2 = 5650
$ws.Cells.Item(24, 11).Value2 = 5650
$ws.Cells.Item(24, 13).Value2 = -5307
$ws.Cells.Item(27, 8).Value2 = 2943.25
$ws.Cells.Item(27, 9).Value2 = 2918.75
$ws.Cells.Item(27, 10).Value2 = 2992.25
$ws.Cells.Item(27, 11).Value2 = 2918.75
$ws.Cells.Item(27, 12).Value2 = 2992.25
$ws.Cells.Item(27, 13).Value2 = -2811.75
$ws.Cells.Item(27, 14).Value2 = -3206.25
$ws.Cells.Item(40, 8).Value2 = 8800.799999999999
$ws.Cells.Item(40, 9).Value2 = 4004
$ws.Cells.Item(40, 11).Value2 = 4004
$ws.Cells.Item(40, 13).Value2 = -3868
$ws.Cells.Item(132, 8).Value2 = 10061
$ws.Cells.Item(132, 9).Value2 = 8335.764999999999
$ws.Cells.Item(132, 11).Value2 = 25007.295
$ws.Cells.Item(132, 13).Value2 = -22477.295
$ws.Cells.Item(136, 8).Value2 = 3349
$ws.Cells.Item(136, 9).Value2 = 3283.4
$ws.Cells.Item(136, 11).Value2 = 9850.200000000001
$ws.Cells.Item(136, 13).Value2 = -7300.200000000001

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(100, 8).Value2 = 865.6667
$ws.Cells.Item(100, 9).Value2 = 865.6667
$ws.Cells.Item(100, 10).Value2 = 0
$ws.Cells.Item(100, 11).Value2 = 1731.3334
$ws.Cells.Item(100, 12).Value2 = 0
$ws.Cells.Item(100, 13).ClearContents()
$ws.Cells.Item(100, 14).Value2 = -1190.3334
$ws.Cells.Item(132, 8).Value2 = 2599.5
$ws.Cells.Item(132, 9).Value2 = 1799.3334
$ws.Cells.Item(132, 10).Value2 = 5000
$ws.Cells.Item(132, 11).Value2 = 5398.0002
$ws.Cells.Item(132, 12).Value2 = 15000
$ws.Cells.Item(132, 13).Value2 = -2868.0002
$ws.Cells.Item(132, 14).Value2 = -20060
$ws.Cells.Item(136, 8).Value2 = 1384.8462
$ws.Cells.Item(136, 9).Value2 = 1384.8462
$ws.Cells.Item(136, 10).Value2 = 0
$ws.Cells.Item(136, 11).Value2 = 4154.5386
$ws.Cells.Item(136, 12).Value2 = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).Value2 = -6750
